$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the default worksheet from the Russian-locale "Лист1" to "Sheet1"
$ws.Name = "Sheet1"

# Move/save the active selection to B2 (was G9)
$ws.Range("B2").Select()
